# Actualizacoes dia 27 as 17
# Applies the cell-level corrections captured in the source diff:
#   - two shared-string text fixes (ID C_1114 -> C_1074, Turma AP-2 -> AP-1)
#   - per-attendee nrSessao / photovoice / apresentacao_photovoice /
#     nivel_engajamento / Familia_Beneficiario corrections across rows 2-48

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N2").Value = 2
$ws.Range("P2").Value = "NAO"
$ws.Range("N3").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("P4").Value = "NAO"
$ws.Range("N5").Value = 6
$ws.Range("P5").Value = "SIM"
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = "SIM"
$ws.Range("P6").Value = "SIM"
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = ""
$ws.Range("N8").Value = 8
$ws.Range("N10").Value = 3
$ws.Range("N11").Value = 5
$ws.Range("N12").Value = 4
$ws.Range("N13").Value = 7
$ws.Range("O13").Value = "SIM"
$ws.Range("P13").Value = "NAO"
$ws.Range("N14").Value = 4
$ws.Range("P14").Value = "NAO"
$ws.Range("N15").Value = 5
$ws.Range("N16").Value = 7
$ws.Range("P16").Value = "SIM"
$ws.Range("N18").Value = 5
$ws.Range("P18").Value = "SIM"
$ws.Range("N19").Value = 4
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = ""
$ws.Range("P20").Value = ""
$ws.Range("N21").Value = 5
$ws.Range("P21").Value = "SIM"
$ws.Range("N22").Value = 8
$ws.Range("O22").Value = ""
$ws.Range("P22").Value = ""
$ws.Range("A23").Value = "C_1074"
$ws.Range("F23").Value = "AP-1"
$ws.Range("L23").Value = "SIM"
$ws.Range("N23").Value = 4
$ws.Range("N24").Value = 6
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = ""
$ws.Range("P25").Value = ""
$ws.Range("N26").Value = 4
$ws.Range("N27").Value = 5
$ws.Range("P27").Value = "NAO"
$ws.Range("N28").Value = 3
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = ""
$ws.Range("P29").Value = ""
$ws.Range("Q29").Value = "NAO ENGAJAD@"
$ws.Range("N30").Value = 2
$ws.Range("O30").Value = "SIM"
$ws.Range("N32").Value = 8
$ws.Range("O32").Value = ""
$ws.Range("P32").Value = ""
$ws.Range("N33").Value = 5
$ws.Range("P33").Value = "SIM"
$ws.Range("Q33").Value = "ENGAJAD@"
$ws.Range("N35").Value = 7
$ws.Range("O35").Value = "NAO"
$ws.Range("N36").Value = 1
$ws.Range("O36").Value = ""
$ws.Range("N37").Value = 7
$ws.Range("P37").Value = "SIM"
$ws.Range("N38").Value = 2
$ws.Range("P38").Value = ""
$ws.Range("N39").Value = 6
$ws.Range("O39").Value = "SIM"
$ws.Range("P39").Value = "NAO"
$ws.Range("Q39").Value = "ENGAJAD@"
$ws.Range("N40").Value = 1
$ws.Range("O40").Value = ""
$ws.Range("P40").Value = ""
$ws.Range("N41").Value = 1
$ws.Range("N42").Value = 2
$ws.Range("P42").Value = "SIM"
$ws.Range("N43").Value = 5
$ws.Range("N44").Value = 2
$ws.Range("O44").Value = "SIM"
$ws.Range("P44").Value = "NAO"
$ws.Range("Q44").Value = "ENGAJAD@"
$ws.Range("N46").Value = 2
$ws.Range("P46").Value = "NAO"
$ws.Range("N47").Value = 6
$ws.Range("O47").Value = "SIM"
$ws.Range("P47").Value = "SIM"
$ws.Range("N48").Value = 8
Write-Output "Updated 89 cell(s) on Sheet1."
